# Daily attendance processing - 2026-01-16 11:07:16
# Reorders the comma-separated "Recorded By" values in column G so that
# "System" (in whatever casing it was originally recorded) moves to the
# front of the list instead of the back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reverse-CommaList($s) {
    if ($s -eq $null) {
        return $s
    }

    $parts = $s.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $n = $trimmed.Length
    if ($n -le 1) {
        return $s
    }

    $result = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $result += $trimmed[$i]
    }

    return ($result -join ", ")
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $newVal = Reverse-CommaList $val

    # Use .Equals() for an ordinal (case-sensitive) comparison, since the
    # default -eq/-ne operators in this engine compare strings
    # case-insensitively and would miss changes that are only a case swap
    # (e.g. "system, ... , System" -> "System, ... , system").
    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
